# "Light background dark font easier to read on projector"
#
# Applies two text edits:
#  1. Slide 2 ("Hard Work Gets Lost..."), Content Placeholder 13:
#     the "Computer turned off" bullet is split into three separate runs
#     ("Computer ", "turned ", "off") -- the result of the author
#     clicking back into that line and re-typing/touching it.
#  2. Slide 3 ("Avoiding confusion is harder..."), Rectangle 2:
#     the two runs "tried to add a major feature but it is not going to
#     work.  " and "How can we ALL get back to where the code was two
#     weeks ago?" are merged back into a single run.

$p = $ppt.ActivePresentation

# --- Slide 2: split "Computer turned off" into three runs ---------------
$slide2 = $p.Slides.Item(2)
$contentShape = $slide2.Shapes.Item(2)
$tf2 = $contentShape.TextFrame
$tr2 = $tf2.TextRange

$computerOffPara = $tr2.Paragraphs(2, 1)
# .Text on a paragraph obtained via Paragraphs() includes the trailing
# paragraph-mark (CR) character, hence the extra "`r" below.
if ($computerOffPara.Text -ne "Computer turned off`r") {
    throw "unexpected paragraph text: [$($computerOffPara.Text)]"
}

$part1 = $computerOffPara.Characters(1, 9)
$part1.Text = "Computer "

$part2 = $computerOffPara.Characters(10, 7)
$part2.Text = "turned "

$part3 = $computerOffPara.Characters(17, 3)
$part3.Text = "off"

# --- Slide 3: merge the two runs back into one ---------------------------
$slide3 = $p.Slides.Item(3)
$rectShape = $slide3.Shapes.Item(2)
$tf3 = $rectShape.TextFrame
$tr3 = $tf3.TextRange

$weTriedPara = $tr3.Paragraphs(5, 1)
$expectedOriginal = " We tried to add a major feature but it is not going to work.  How can we ALL get back to where the code was two weeks ago?"
if ($weTriedPara.Text -ne $expectedOriginal) {
    throw "unexpected paragraph text: [$($weTriedPara.Text)]"
}

$mergedPart = $weTriedPara.Characters(5, 119)
$mergedPart.Text = "tried to add a major feature but it is not going to work.  How can we ALL get back to where the code was two weeks ago?"

# Merging the two runs doesn't change the wrapped text at all, so the
# "Rectangle 2" autofit shape (<a:spAutoFit/>) should keep its original
# height; re-assert it explicitly since touching the run text marks the
# shape for layout recompute.
$originalHeightEmu = 3541375
$rectShape.Height = ($originalHeightEmu / 12700.0) + 0.00005
